$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 7.453102666666666
$ws.Range("H2").Value = 22.359308
$ws.Range("I2").Value = 0.0465828813588002
$ws.Range("J2").Value = 0.0465828813588002
$ws.Range("M2").Value = 3.425446666666666
$ws.Range("N2").Value = 10.27634
$ws.Range("O2").Value = 0.6657953389778073
$ws.Range("P2").Value = 0.6657953389778073
$ws.Range("Q2").Value = 25.53020568585777
$ws.Range("R2").Value = 229.77185117272
$ws.Range("S2").Value = 0.03101466528484536
$ws.Range("T2").Value = 0.03101466528484536
$ws.Range("G3").Value = 7.453102666666666
$ws.Range("H3").Value = 22.359308
$ws.Range("I3").Value = 0.0465828813588002
$ws.Range("J3").Value = 0.0465828813588002
$ws.Range("O3").Value = 0.2094791321596951
$ws.Range("P3").Value = 0.2094791321596952
$ws.Range("Q3").Value = 8.032566492794667
$ws.Range("R3").Value = 72.29309843515199
$ws.Range("S3").Value = 0.009758141560539505
$ws.Range("T3").Value = 0.009758141560539507
$ws.Range("G4").Value = 7.453102666666666
$ws.Range("H4").Value = 22.359308
$ws.Range("I4").Value = 0.0465828813588002
$ws.Range("J4").Value = 0.0465828813588002
$ws.Range("M4").Value = 0.62317
$ws.Range("N4").Value = 1.86951
$ws.Range("O4").Value = 0.1211239647746572
$ws.Range("P4").Value = 0.1211239647746572
$ws.Range("Q4").Value = 4.644549988786666
$ws.Range("R4").Value = 41.80094989908
$ws.Range("S4").Value = 0.00564230328080535
$ws.Range("T4").Value = 0.00564230328080535
$ws.Range("G5").Value = 7.453102666666666
$ws.Range("H5").Value = 22.359308
$ws.Range("I5").Value = 0.0465828813588002
$ws.Range("J5").Value = 0.0465828813588002
$ws.Range("M5").Value = 0.01852966666666667
$ws.Range("N5").Value = 0.055589
$ws.Range("O5").Value = 0.003601564087840353
$ws.Range("P5").Value = 0.003601564087840353
$ws.Range("Q5").Value = 0.1381035080457778
$ws.Range("R5").Value = 1.242931572412
$ws.Range("S5").Value = 0.0001677712326099826
$ws.Range("T5").Value = 0.0001677712326099826
$ws.Range("I6").Value = 0.7985586056263754
$ws.Range("J6").Value = 0.7985586056263754
$ws.Range("M6").Value = 3.425446666666666
$ws.Range("N6").Value = 10.27634
$ws.Range("O6").Value = 0.6657953389778073
$ws.Range("P6").Value = 0.6657953389778073
$ws.Range("Q6").Value = 437.6578876008422
$ws.Range("R6").Value = 3938.920988407579
$ws.Range("S6").Value = 0.5316765975266577
$ws.Range("T6").Value = 0.5316765975266577
$ws.Range("I7").Value = 0.7985586056263754
$ws.Range("J7").Value = 0.7985586056263754
$ws.Range("O7").Value = 0.2094791321596951
$ws.Range("P7").Value = 0.2094791321596952
$ws.Range("S7").Value = 0.1672813636852694
$ws.Range("T7").Value = 0.1672813636852694
$ws.Range("I8").Value = 0.7985586056263754
$ws.Range("J8").Value = 0.7985586056263754
$ws.Range("M8").Value = 0.62317
$ws.Range("N8").Value = 1.86951
$ws.Range("O8").Value = 0.1211239647746572
$ws.Range("P8").Value = 0.1211239647746572
$ws.Range("Q8").Value = 79.62035096626333
$ws.Range("R8").Value = 716.58315869637
$ws.Range("S8").Value = 0.09672458441838845
$ws.Range("T8").Value = 0.09672458441838845
$ws.Range("I9").Value = 0.7985586056263754
$ws.Range("J9").Value = 0.7985586056263754
$ws.Range("M9").Value = 0.01852966666666667
$ws.Range("N9").Value = 0.055589
$ws.Range("O9").Value = 0.003601564087840353
$ws.Range("P9").Value = 0.003601564087840353
$ws.Range("Q9").Value = 2.367473664149222
$ws.Range("R9").Value = 21.307262977343
$ws.Range("S9").Value = 0.002876059996059821
$ws.Range("T9").Value = 0.002876059996059821
$ws.Range("G10").Value = 24.73806933333333
$ws.Range("H10").Value = 74.214208
$ws.Range("I10").Value = 0.1546162182837376
$ws.Range("J10").Value = 0.1546162182837376
$ws.Range("M10").Value = 3.425446666666666
$ws.Range("N10").Value = 10.27634
$ws.Range("O10").Value = 0.6657953389778073
$ws.Range("P10").Value = 0.6657953389778073
$ws.Range("Q10").Value = 84.73893713763555
$ws.Range("R10").Value = 762.65043423872
$ws.Range("S10").Value = 0.1029427574636877
$ws.Range("T10").Value = 0.1029427574636877
$ws.Range("G11").Value = 24.73806933333333
$ws.Range("H11").Value = 74.214208
$ws.Range("I11").Value = 0.1546162182837376
$ws.Range("J11").Value = 0.1546162182837376
$ws.Range("O11").Value = 0.2094791321596951
$ws.Range("P11").Value = 0.2094791321596952
$ws.Range("Q11").Value = 26.66140474786133
$ws.Range("R11").Value = 239.952642730752
$ws.Range("S11").Value = 0.03238887122389134
$ws.Range("T11").Value = 0.03238887122389135
$ws.Range("G12").Value = 24.73806933333333
$ws.Range("H12").Value = 74.214208
$ws.Range("I12").Value = 0.1546162182837376
$ws.Range("J12").Value = 0.1546162182837376
$ws.Range("M12").Value = 0.62317
$ws.Range("N12").Value = 1.86951
$ws.Range("O12").Value = 0.1211239647746572
$ws.Range("P12").Value = 0.1211239647746572
$ws.Range("Q12").Value = 15.41602266645333
$ws.Range("R12").Value = 138.74420399808
$ws.Range("S12").Value = 0.01872772937699014
$ws.Range("T12").Value = 0.01872772937699014
$ws.Range("G13").Value = 24.73806933333333
$ws.Range("H13").Value = 74.214208
$ws.Range("I13").Value = 0.1546162182837376
$ws.Range("J13").Value = 0.1546162182837376
$ws.Range("M13").Value = 0.01852966666666667
$ws.Range("N13").Value = 0.055589
$ws.Range("O13").Value = 0.003601564087840353
$ws.Range("P13").Value = 0.003601564087840353
$ws.Range("Q13").Value = 0.4583881787235555
$ws.Range("R13").Value = 4.125493608512
$ws.Range("S13").Value = 0.0005568602191683943
$ws.Range("T13").Value = 0.0005568602191683943
$ws.Range("G14").Value = 0.03876633333333333
$ws.Range("H14").Value = 0.116299
$ws.Range("I14").Value = 0.0002422947310868075
$ws.Range("J14").Value = 0.0002422947310868075
$ws.Range("M14").Value = 3.425446666666666
$ws.Range("N14").Value = 10.27634
$ws.Range("O14").Value = 0.6657953389778073
$ws.Range("P14").Value = 0.6657953389778073
$ws.Range("Q14").Value = 0.1327920072955555
$ws.Range("R14").Value = 1.19512806566
$ws.Range("S14").Value = 0.0001613187026164777
$ws.Range("T14").Value = 0.0001613187026164777
$ws.Range("G15").Value = 0.03876633333333333
$ws.Range("H15").Value = 0.116299
$ws.Range("I15").Value = 0.0002422947310868075
$ws.Range("J15").Value = 0.0002422947310868075
$ws.Range("O15").Value = 0.2094791321596951
$ws.Range("P15").Value = 0.2094791321596952
$ws.Range("Q15").Value = 0.04178033821733333
$ws.Range("R15").Value = 0.376023043956
$ws.Range("S15").Value = 0.00005075568999493115009960173
$ws.Range("T15").Value = 0.000050755689994931156875865308
$ws.Range("G16").Value = 0.03876633333333333
$ws.Range("H16").Value = 0.116299
$ws.Range("I16").Value = 0.0002422947310868075
$ws.Range("J16").Value = 0.0002422947310868075
$ws.Range("M16").Value = 0.62317
$ws.Range("N16").Value = 1.86951
$ws.Range("O16").Value = 0.1211239647746572
$ws.Range("P16").Value = 0.1211239647746572
$ws.Range("Q16").Value = 0.02415801594333333
$ws.Range("R16").Value = 0.21742214349
$ws.Range("S16").Value = 0.000029347698473243509373559537
$ws.Range("T16").Value = 0.000029347698473243509373559537
$ws.Range("G17").Value = 0.03876633333333333
$ws.Range("H17").Value = 0.116299
$ws.Range("I17").Value = 0.0002422947310868075
$ws.Range("J17").Value = 0.0002422947310868075
$ws.Range("M17").Value = 0.01852966666666667
$ws.Range("N17").Value = 0.055589
$ws.Range("O17").Value = 0.003601564087840353
$ws.Range("P17").Value = 0.003601564087840353
$ws.Range("Q17").Value = 0.0007183272345555556
$ws.Range("R17").Value = 0.006464945111
$ws.Range("S17").Value = 0.000000872640002155181589069664
$ws.Range("T17").Value = 0.000000872640002155181589069664
